$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = [double]"-414685755453185152"
$ws.Range("C2").Value = [double]"-172048015514311"
$ws.Range("D2").Value = [double]"101689501379444527333376"
$ws.Range("E2").Value = [double]"9834185715030458368"
$ws.Range("B3").Value = [double]"-414103148410964800"
$ws.Range("C3").Value = [double]"-171805981845942"
$ws.Range("D3").Value = [double]"101552918427403292770304"
$ws.Range("E3").Value = [double]"9820316061320046592"
$ws.Range("B4").Value = [double]"-413520581335863296"
$ws.Range("C4").Value = [double]"-171563945109500"
$ws.Range("D4").Value = [double]"101416381568627412828160"
$ws.Range("E4").Value = [double]"9806446419267573760"
$ws.Range("B5").Value = [double]"-412938045856773760"
$ws.Range("C5").Value = [double]"-171321463870520"
$ws.Range("D5").Value = [double]"101279851876151065575424"
$ws.Range("E5").Value = [double]"9792576787378438144"
$ws.Range("B6").Value = [double]"-412355541057548224"
$ws.Range("C6").Value = [double]"-171078398708179"
$ws.Range("D6").Value = [double]"101143351766159684796416"
$ws.Range("E6").Value = [double]"9778707166723219456"
$ws.Range("B7").Value = [double]"-411773040573451264"
$ws.Range("C7").Value = [double]"-170834977968150"
$ws.Range("D7").Value = [double]"101006886204998412140544"
$ws.Range("E7").Value = [double]"9764837546392616960"
$ws.Range("B8").Value = [double]"-411190527005611648"
$ws.Range("C8").Value = [double]"-170591592531619"
$ws.Range("D8").Value = [double]"100870207038307352707072"
$ws.Range("E8").Value = [double]"9750967922969516032"
$ws.Range("B9").Value = [double]"-410608010284447808"
$ws.Range("C9").Value = [double]"-170350024380500"
$ws.Range("D9").Value = [double]"100733437300810723622912"
$ws.Range("E9").Value = [double]"9737098307609286656"
$ws.Range("B10").Value = [double]"-410025538796957696"
$ws.Range("C10").Value = [double]"-170107192542451"
$ws.Range("D10").Value = [double]"100596939768077198819328"
$ws.Range("E10").Value = [double]"9723228697136887808"
$ws.Range("B11").Value = [double]"-409443145894310528"
$ws.Range("C11").Value = [double]"-169862203101715"
$ws.Range("D11").Value = [double]"100460307232916055785472"
$ws.Range("E11").Value = [double]"9709359113879023616"
$ws.Range("B12").Value = [double]"-408860697143694656"
$ws.Range("C12").Value = [double]"-169618426336256"
$ws.Range("D12").Value = [double]"100323684236277323923456"
$ws.Range("E12").Value = [double]"9695489511912505344"
$ws.Range("B13").Value = [double]"-408278261450843264"
$ws.Range("C13").Value = [double]"-169374883766055"
$ws.Range("D13").Value = [double]"100187138660363090788352"
$ws.Range("E13").Value = [double]"9681619910661543936"
$ws.Range("B14").Value = [double]"-407695808727930816"
$ws.Range("C14").Value = [double]"-169130413702686"
$ws.Range("D14").Value = [double]"100050871864227541286912"
$ws.Range("E14").Value = [double]"9667750307175444480"
$ws.Range("B15").Value = [double]"-407113339877007488"
$ws.Range("C15").Value = [double]"-168882902777658"
$ws.Range("D15").Value = [double]"99914434130193806262272"
$ws.Range("E15").Value = [double]"9653880698307862528"
$ws.Range("B16").Value = [double]"-406530935205805696"
$ws.Range("C16").Value = [double]"-168635123179448"
$ws.Range("D16").Value = [double]"99778069783030146269184"
$ws.Range("E16").Value = [double]"9640011110852474880"
$ws.Range("B17").Value = [double]"-405948490449232448"
$ws.Range("C17").Value = [double]"-168385134099680"
$ws.Range("D17").Value = [double]"99641636674204011069440"
$ws.Range("E17").Value = [double]"9626141511763718144"
$ws.Range("B18").Value = [double]"-405365966137969280"
$ws.Range("C18").Value = [double]"-168135670586875"
$ws.Range("D18").Value = [double]"99505291342404393107456"
$ws.Range("E18").Value = [double]"9612271884470534144"
$ws.Range("B19").Value = [double]"-404783567166120512"
$ws.Range("C19").Value = [double]"-167889231429461"
$ws.Range("D19").Value = [double]"99369010004353037107200"
$ws.Range("E19").Value = [double]"9598402297237733376"
$ws.Range("B20").Value = [double]"-404201105560723392"
$ws.Range("C20").Value = [double]"-167647101994698"
$ws.Range("D20").Value = [double]"99232636886439581188096"
$ws.Range("E20").Value = [double]"9584532688124192768"
$ws.Range("B21").Value = [double]"-403618766980988800"
$ws.Range("C21").Value = [double]"-167399480395769"
$ws.Range("D21").Value = [double]"99096348829084645588992"
$ws.Range("E21").Value = [double]"9570663122144473088"
$ws.Range("B22").Value = [double]"-403036279618939072"
$ws.Range("C22").Value = [double]"-167151438395178"
$ws.Range("D22").Value = [double]"98959875858994510692352"
$ws.Range("E22").Value = [double]"9556793512277878784"
$ws.Range("B23").Value = [double]"-402453962995075456"
$ws.Range("C23").Value = [double]"-166905107643026"
$ws.Range("D23").Value = [double]"98823958405789708713984"
$ws.Range("E23").Value = [double]"9542923951324366848"
$ws.Range("B24").Value = [double]"-401871607086526720"
$ws.Range("C24").Value = [double]"-166659364068245"
$ws.Range("D24").Value = [double]"98687937481902978498560"
$ws.Range("E24").Value = [double]"9529054378289557504"
$ws.Range("B25").Value = [double]"-401289229841510400"
$ws.Range("C25").Value = [double]"-166414578469918"
$ws.Range("D25").Value = [double]"98551684711607682727936"
$ws.Range("E25").Value = [double]"9515184800921980928"
$ws.Range("B26").Value = [double]"-400706863195345408"
$ws.Range("C26").Value = [double]"-166172141070565"
$ws.Range("D26").Value = [double]"98415446032347945762816"
$ws.Range("E26").Value = [double]"9501315223276883968"
$ws.Range("B27").Value = [double]"-400124451776021760"
$ws.Range("C27").Value = [double]"-165927013953695"
$ws.Range("D27").Value = [double]"98279397793893934170112"
$ws.Range("E27").Value = [double]"9487445630468073472"
$ws.Range("B28").Value = [double]"-399542051159133952"
$ws.Range("C28").Value = [double]"-165679984863616"
$ws.Range("D28").Value = [double]"98143285890155106271232"
$ws.Range("E28").Value = [double]"9473576036040912896"
